$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.619.60'
$ws.Range('E2').Value = '  +1.37%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.881.02'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7275'
$ws.Range('E5').Value = '  +3.32%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '239.84'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9994'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07906'
$ws.Range('E8').Value = '  -3.52%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3095'
$ws.Range('E9').Value = '  +2.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '25.21'
$ws.Range('E10').Value = '  +8.22%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08250'
$ws.Range('E11').Value = '  +1.09%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.889.16'
$ws.Range('E12').Value = '  +1.80%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.7281'
$ws.Range('E13').Value = '  +2.73%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.271'
$ws.Range('E14').Value = '  +2.02%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '90.57'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '29.656.79'
$ws.Range('E16').Value = '  +1.49%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.871'
$ws.Range('E17').Value = '  +1.53%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000007899'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '243.15'
$ws.Range('E19').Value = '  +2.97%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.42'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.141.28'
$ws.Range('E21').Value = '  +2.58%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.0000'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '7.786'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1605'
$ws.Range('E25').Value = '  +11.44%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '163.02'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.006'
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.43'
$ws.Range('E28').Value = '  +1.97%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.955'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.372'
$ws.Range('E30').Value = '  -3.81%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.483'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.364'
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.116'
$ws.Range('E33').Value = '  +1.41%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05278'
$ws.Range('E34').Value = '  +1.41%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.200'
$ws.Range('E35').Value = '  +2.62%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7198'
$ws.Range('E36').Value = '  +1.84%  '
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.667'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01870'
$ws.Range('E39').Value = '  +1.12%  '
$ws.Range('E40').Value = '  -0.45%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.191.04'
$ws.Range('E41').Value = '  +4.27%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.9035'
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.018'
$ws.Range('E43').Value = '  +2.52%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '72.04'
$ws.Range('E44').Value = '  +2.57%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.4325'
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.9996'
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '102.89'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5371'
$ws.Range('E48').Value = '  -0.95%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.784'
$ws.Range('E49').Value = '  +0.52%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.258'
$ws.Range('E50').Value = '  +0.61%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.886'
$ws.Range('E51').Value = '  +5.15%  '
